$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chrom to dia")
$ws.Activate()

# Header for column B
$ws.Range("B1").Value = "cc % 12"
$ws.Range("B1").NumberFormat = "@"

# Formula for B3:B42 -> MOD(A,12), entered first so Excel groups them as one
# shared-formula block; B2 entered afterwards stays an ungrouped formula.
$ws.Range("B3:B42").Formula = "=MOD(A3, 12)"
$ws.Range("B2").Formula = "=MOD(A2, 12)"

# Column A autofit-like width (bestFit) to match diff (6.44140625)
$ws.Columns.Item(1).ColumnWidth = 5.6

# Selection as in diff
$ws.Range("B16").Select()

$wb.Save()
